# Danh sách khách hàng còn dư nợ tại HỆ THỐNG
# fixbug tinh chiết khấu đơn thu nợ
#
# A new customer ("Cầm dương", mã KH 443) needs to be inserted at the top
# of the data list (row 2, right under the header row), pushing every
# existing customer row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 2 - shifts all existing customer rows down by one.
$ws.Rows("2:2").Insert()

# Fill in the new customer's data in the freshly inserted row.
$ws.Range("A2").Value = "KH"
$ws.Range("B2").Value = 443
$ws.Range("C2").Value = "Cầm dương"
$ws.Range("D2").Value = "LONG XUYÊN"
$ws.Range("E2").Value = 0
# Leading zero must be kept literal, so force the cell to text ('-prefix),
# then drop back to the Normal style so no stray text format lingers.
$ws.Range("F2").Value = "'0364580162"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 2000000
$ws.Range("J2").Value = 12000000
